$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The inventory rows 31-36 got re-ordered upstream (row 36 "Chip Epson / C9345"
# moved up to row 31, pushing the previous rows 31-35 down by one). Apply the
# new row contents (columns A-G and J; H and I are formulas and recompute
# automatically from D/E/F/G).

$rows = @(
    @{ Row = 31; A = "HGP1Q2"; B = "Chip Epson"; C = "C9345"; D = 0; E = 100000; F = 20; G = 0; J = 0 },
    @{ Row = 32; A = "LJNL6J"; B = "Drum DL-410 para Tambor de Unidad de imagen Pantum"; C = "P3010DW P3300DN P3300DW M6700DW M6800FDW M7100DN M7200FDW M7300FDW M7300FDN"; D = 0; E = 100000; F = 9; G = 0; J = 0 },
    @{ Row = 33; A = "UAYYDQ"; B = "DMD para proyector 8060-6039B"; C = "Benq MP515 MP515ST NEC NP115 OPTOMA ES526"; D = 200000; E = 400000; F = 1; G = 1; J = 200000 },
    @{ Row = 34; A = "CJVMIV"; B = "Correa de transporte Epson"; C = "TM U950"; D = 0; E = 100000; F = 7; G = 3; J = 0 },
    @{ Row = 35; A = "Y6I8Q7"; B = "Correa de plotter  HP DesignJet"; C = "500 510 800 815 de 24`""; D = 0; E = 350000; F = 3; G = 0; J = 0 },
    @{ Row = 36; A = "O3F6AE"; B = "Correa de plotter HP DesignJet"; C = "500 510 800 815 de 42`""; D = 0; E = 350000; F = 2; G = 0; J = 0 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value2 = $r.A
    $ws.Cells.Item($n, 2).Value2 = $r.B
    $ws.Cells.Item($n, 3).Value2 = $r.C
    $ws.Cells.Item($n, 4).Value2 = $r.D
    $ws.Cells.Item($n, 5).Value2 = $r.E
    $ws.Cells.Item($n, 6).Value2 = $r.F
    $ws.Cells.Item($n, 7).Value2 = $r.G
    $ws.Cells.Item($n, 10).Value2 = $r.J
}
